$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [DateTime]"2021-05-27"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 47
$ws.Range("R2").Value = "Región Metropolitana"

# Row 3
$ws.Range("D3").Value = [DateTime]"2021-05-27"
$ws.Range("M3").Value = 50
$ws.Range("R3").Value = "Región Metropolitana"

# Row 4
$ws.Range("D4").Value = [DateTime]"2021-05-27"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 800

# Row 5
$ws.Range("D5").Value = [DateTime]"2021-05-03"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 68
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 57
$ws.Range("N6").Value = 8000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 8000
$ws.Range("S6").Value = 800

# Row 7
$ws.Range("D7").Value = [DateTime]"2021-04-28"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 47
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 9000
$ws.Range("S7").Value = 900

# Row 8
$ws.Range("D8").Value = [DateTime]"2021-05-12"
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 800

# Row 9
$ws.Range("D9").Value = [DateTime]"2021-05-12"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 7000
$ws.Range("S9").Value = 700

# Row 10
$ws.Range("D10").Value = [DateTime]"2021-05-17"
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 58
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("D11").Value = [DateTime]"2021-05-17"
$ws.Range("M11").Value = 65

# Row 12
$ws.Range("D12").Value = [DateTime]"2021-05-17"
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("S12").Value = 800

# Row 13
$ws.Range("D13").Value = [DateTime]"2021-05-10"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 1000

# Row 14
$ws.Range("D14").Value = [DateTime]"2021-05-10"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 67
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("S14").Value = 800

# Row 15
$ws.Range("D15").Value = [DateTime]"2021-04-16"
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("R15").Value = "Provincia de Quillota"
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = [DateTime]"2021-04-22"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 45
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("R16").Value = "Provincia de Quillota"
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = [DateTime]"2021-04-22"
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 48
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 800

# Row 18
$ws.Range("D18").Value = [DateTime]"2021-05-05"
$ws.Range("M18").Value = 58

# Row 19
$ws.Range("D19").Value = [DateTime]"2021-04-20"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 10000
$ws.Range("S19").Value = 1000

# Row 22
$ws.Range("D22").Value = [DateTime]"2021-04-15"

# Row 23
$ws.Range("D23").Value = [DateTime]"2021-04-26"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 48
$ws.Range("R23").Value = "Provincia de Quillota"

# Row 24
$ws.Range("D24").Value = [DateTime]"2021-05-13"
$ws.Range("M24").Value = 56

# Row 25
$ws.Range("D25").Value = [DateTime]"2021-05-13"
$ws.Range("M25").Value = 50

# Row 26
$ws.Range("D26").Value = [DateTime]"2021-04-29"
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("S26").Value = 1000

# Row 27
$ws.Range("D27").Value = [DateTime]"2021-05-06"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 56
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("S27").Value = 1000

# Row 28
$ws.Range("D28").Value = [DateTime]"2021-05-06"
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("S28").Value = 800

# Row 30
$ws.Range("D30").Value = [DateTime]"2021-04-23"
